# Auto-generated Excel COM-interop script
# Applies the scheduled price/profit refresh to the 8 leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 837.8
$ws.Range("I18").Value = 700
$ws.Range("K18").Value = 700
$ws.Range("M18").Value = -416
$ws.Range("H19").Value = 556.2
$ws.Range("I19").Value = 656.75
$ws.Range("K19").Value = 656.75
$ws.Range("M19").Value = -481.75
$ws.Range("H64").Value = 8594.299999999999
$ws.Range("I64").Value = 7331.6665
$ws.Range("J64").Value = 9135.429
$ws.Range("K64").Value = 7331.6665
$ws.Range("L64").Value = 9135.429
$ws.Range("M64").Value = -7083.6665
$ws.Range("N64").Value = -9631.429
$ws.Range("H67").Value = 8594.299999999999
$ws.Range("I67").Value = 7331.6665
$ws.Range("J67").Value = 9135.429
$ws.Range("K67").Value = 7331.6665
$ws.Range("L67").Value = 9135.429
$ws.Range("M67").Value = -6473.6665
$ws.Range("N67").Value = -10851.429
$ws.Range("H92").Value = 2144
$ws.Range("I92").Value = 1128.5
$ws.Range("K92").Value = 1128.5
$ws.Range("M92").Value = 119.5
$ws.Range("H96").Value = 290.92307
$ws.Range("I96").Value = 193.81818
$ws.Range("K96").Value = 581.4545400000001
$ws.Range("M96").Value = 791.5454599999999
$ws.Range("H112").Value = 2204.5
$ws.Range("I112").Value = 2000
$ws.Range("J112").Value = 2272.6667
$ws.Range("K112").Value = 6000
$ws.Range("L112").Value = 6818.000100000001
$ws.Range("M112").Value = -4892
$ws.Range("N112").Value = -9034.000100000001
$ws.Range("H127").Value = 897.8
$ws.Range("I127").Value = 897.8
$ws.Range("K127").Value = 2693.4
$ws.Range("M127").Value = 2266.6
$ws.Range("H135").Value = 2190.8572
$ws.Range("I135").Value = 2000.75
$ws.Range("K135").Value = 18006.75
$ws.Range("M135").Value = -15471.75
$ws.Range("H137").Value = 2241
$ws.Range("I137").Value = 1530.6666
$ws.Range("J137").Value = 4016.8333
$ws.Range("K137").Value = 4591.9998
$ws.Range("L137").Value = 12050.4999
$ws.Range("M137").Value = -2041.9998
$ws.Range("N137").Value = -17150.4999
$ws.Range("H141").Value = 2056
$ws.Range("I141").Value = 1059.8
$ws.Range("J141").Value = 3716.3333
$ws.Range("K141").Value = 3179.4
$ws.Range("L141").Value = 11148.9999
$ws.Range("M141").Value = 2000.6
$ws.Range("N141").Value = -21508.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1109.8
$ws.Range("I97").Value = 657
$ws.Range("J97").Value = 2355
$ws.Range("K97").Value = 657
$ws.Range("L97").Value = 2355
$ws.Range("M97").Value = -161
$ws.Range("N97").Value = -3347
$ws.Range("H128").Value = 175000
$ws.Range("J128").Value = 175000
$ws.Range("L128").Value = 175000
$ws.Range("N128").Value = -184960

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4531.1113
$ws.Range("I20").Value = 4694.4
$ws.Range("J20").Value = 4327
$ws.Range("K20").Value = 4694.4
$ws.Range("L20").Value = 4327
$ws.Range("M20").Value = -4447.4
$ws.Range("N20").Value = -4821
$ws.Range("H80").Value = 171.5
$ws.Range("I80").Value = 114.333336
$ws.Range("J80").Value = 190.55556
$ws.Range("K80").Value = 114.333336
$ws.Range("L80").Value = 190.55556
$ws.Range("M80").Value = 883.666664
$ws.Range("N80").Value = -2186.55556
$ws.Range("H83").Value = 171.5
$ws.Range("I83").Value = 114.333336
$ws.Range("J83").Value = 190.55556
$ws.Range("K83").Value = 571.66668
$ws.Range("L83").Value = 952.7778000000001
$ws.Range("M83").Value = 4420.33332
$ws.Range("N83").Value = -10936.7778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 500783.8
$ws.Range("J41").Value = 541862.2
$ws.Range("L41").Value = 541862.2
$ws.Range("N41").Value = -542718.2
$ws.Range("H58").Value = 2677.3572
$ws.Range("I58").Value = 954.4545000000001
$ws.Range("K58").Value = 954.4545000000001
$ws.Range("M58").Value = -751.4545000000001
$ws.Range("H99").Value = 4367
$ws.Range("I99").Value = 3767.4614
$ws.Range("K99").Value = 3767.4614
$ws.Range("M99").Value = -2269.4614
$ws.Range("H126").Value = 4367
$ws.Range("I126").Value = 3767.4614
$ws.Range("K126").Value = 11302.3842
$ws.Range("M126").Value = -8832.3842
$ws.Range("H132").Value = 4901.5713
$ws.Range("J132").Value = 6387.3335
$ws.Range("L132").Value = 19162.0005
$ws.Range("N132").Value = -24222.0005
$ws.Range("H136").Value = 2677.3572
$ws.Range("I136").Value = 954.4545000000001
$ws.Range("K136").Value = 2863.3635
$ws.Range("M136").Value = -313.3635000000004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 334007.34
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 13401200
$ws.Range("J11").Value = 8003500
$ws.Range("L11").Value = 8003500
$ws.Range("N11").Value = -8003778
$ws.Range("H24").Value = 44500
$ws.Range("J24").Value = 44500
$ws.Range("L24").Value = 44500
$ws.Range("N24").Value = -44846
$ws.Range("H62").Value = 28995
$ws.Range("I62").Value = 28995
$ws.Range("K62").Value = 28995
$ws.Range("M62").Value = -28309
$ws.Range("H65").Value = 28995
$ws.Range("I65").Value = 28995
$ws.Range("K65").Value = 86985
$ws.Range("M65").Value = -83553
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H70").Value = 5376.091
$ws.Range("I70").Value = 4491.4287
$ws.Range("J70").Value = 6924.25
$ws.Range("K70").Value = 4491.4287
$ws.Range("L70").Value = 6924.25
$ws.Range("M70").Value = -4221.4287
$ws.Range("N70").Value = -7464.25
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H73").Value = 5376.091
$ws.Range("I73").Value = 4491.4287
$ws.Range("J73").Value = 6924.25
$ws.Range("K73").Value = 4491.4287
$ws.Range("L73").Value = 6924.25
$ws.Range("M73").Value = -3555.4287
$ws.Range("N73").Value = -8796.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1649.5
$ws.Range("I22").Value = 1779.4
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 1779.4
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -1484.4
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 1649.5
$ws.Range("I27").Value = 1779.4
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 1779.4
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -1672.4
$ws.Range("N27").Value = -1214
$ws.Range("H40").Value = 6997.6665
$ws.Range("I40").Value = 5664.6665
$ws.Range("J40").Value = 8330.666999999999
$ws.Range("K40").Value = 5664.6665
$ws.Range("L40").Value = 8330.666999999999
$ws.Range("M40").Value = -5528.6665
$ws.Range("N40").Value = -8602.666999999999
$ws.Range("H46").Value = 6179.95
$ws.Range("I46").Value = 5639.125
$ws.Range("J46").Value = 6540.5
$ws.Range("K46").Value = 5639.125
$ws.Range("L46").Value = 6540.5
$ws.Range("M46").Value = -5451.125
$ws.Range("N46").Value = -6916.5
$ws.Range("H55").Value = 805.17645
$ws.Range("I55").Value = 934.6667
$ws.Range("J55").Value = 494.4
$ws.Range("K55").Value = 934.6667
$ws.Range("L55").Value = 494.4
$ws.Range("M55").Value = -761.6667
$ws.Range("N55").Value = -840.4
$ws.Range("H132").Value = 4999
$ws.Range("I132").Value = 4999
$ws.Range("K132").Value = 14997
$ws.Range("M132").Value = -12467
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 9999
$ws.Range("J31").Value = 9999
$ws.Range("L31").Value = 9999
$ws.Range("N31").Value = -10695
$ws.Range("H81").Value = 3196.4
$ws.Range("J81").Value = 3196.4
$ws.Range("L81").Value = 6392.8
$ws.Range("N81").Value = -8514.799999999999
$ws.Range("H84").Value = 3196.4
$ws.Range("J84").Value = 3196.4
$ws.Range("L84").Value = 31964
$ws.Range("N84").Value = -42572
$ws.Range("H122").Value = 3183.1667
$ws.Range("I122").Value = 1734.3334
$ws.Range("K122").Value = 5203.0002
$ws.Range("M122").Value = -2753.0002
